$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "partnership" sheet - replace the 1994-2023 series with the new
#    2010-2023 series (rates evaluated from the new partnership alignment)
#    and drop the (unused) cell style that used to be applied to col A.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("partnership")

# wipe out everything below the header row (content + formatting) so the
# leftover rows (16:31) disappear completely instead of lingering as blank
# styled cells
$ws.Range("A2:B31").Clear()

$years = 2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020,2021,2022,2023
$shares = 0.5763,0.5694,0.5716,0.5646,0.5563,0.546,0.5481,0.541,0.5513,0.5539,0.5525,0.5526,0.5527,0.5528

for ($i = 0; $i -lt $years.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $years[$i]
    $ws.Cells.Item($r, 2).Value2 = $shares[$i]
}

# ---------------------------------------------------------------------------
# 2. "raw data" sheet - the W:Z column block used to carry a dedicated
#    (blank) style; the partnership-alignment rework drops it, which also
#    clears the stray style="3" column formatting + per-cell style refs
#    that sat on top of the W/X/Y helper columns (used to build the
#    partnership tab above).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("raw data")

$ws2.Range("X3:Y3").Style = "Normal"
$ws2.Range("W4:Y33").Style = "Normal"
$ws2.Range("W1:Z1").EntireColumn.ClearFormats()
